# Apply the crypto price/volume/hour-counter refresh described in the commit
# message ("Updated symbol list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.26"
$ws.Range("E2").Value = "'-0.68%"
$ws.Range("G2").Value = "'21"
$ws.Range("D3").Value = "'27.16"
$ws.Range("E3").Value = "'3.77%"
$ws.Range("G3").Value = "'21"
$ws.Range("D4").Value = "'5.118"
$ws.Range("E4").Value = "'0.74%"
$ws.Range("G4").Value = "'21"
$ws.Range("E5").Value = "'1.67%"
$ws.Range("G5").Value = "'21"
$ws.Range("D6").Value = "'6.515"
$ws.Range("E6").Value = "'0.66%"
$ws.Range("G6").Value = "'21"
$ws.Range("D7").Value = "'0.8192"
$ws.Range("E7").Value = "'0.77%"
$ws.Range("G7").Value = "'21"
$ws.Range("D8").Value = "'0.8589"
$ws.Range("E8").Value = "'1.80%"
$ws.Range("G8").Value = "'21"
$ws.Range("B9").Value = "MandalaExchangeToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D9").Value = "'0.06937"
$ws.Range("E9").Value = "'-0.94%"
$ws.Range("G9").Value = "'21"
$ws.Range("B10").Value = "BitrueCoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D10").Value = "'0.02857"
$ws.Range("E10").Value = "'0.71%"
$ws.Range("G10").Value = "'21"
$ws.Range("B11").Value = "BitMartToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D11").Value = "'0.09391"
$ws.Range("E11").Value = "'-0.13%"
$ws.Range("G11").Value = "'21"
$ws.Range("B12").Value = "BitForexToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D12").Value = "'0.001515"
$ws.Range("E12").Value = "'-1.35%"
$ws.Range("G12").Value = "'21"
$ws.Range("B13").Value = "CoinExToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D13").Value = "'0.04052"
$ws.Range("E13").Value = "'-12.97%"
$ws.Range("G13").Value = "'21"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006028"
$ws.Range("E14").Value = "'0.30%"
$ws.Range("G14").Value = "'21"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006213"
$ws.Range("E15").Value = "'0.96%"
$ws.Range("G15").Value = "'21"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.509"
$ws.Range("E16").Value = "'-2.72%"
$ws.Range("G16").Value = "'21"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.012"
$ws.Range("E17").Value = "'-0.02%"
$ws.Range("G17").Value = "'21"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.318"
$ws.Range("E18").Value = "'12.79%"
$ws.Range("G18").Value = "'21"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3165"
$ws.Range("E19").Value = "'1.23%"
$ws.Range("G19").Value = "'21"
$ws.Range("B20").Value = "WazirX"
$ws.Range("C20").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D20").Value = "'0.1331"
$ws.Range("E20").Value = "'-0.35%"
$ws.Range("G20").Value = "'21"
$ws.Range("D21").Value = "'0.03226"
$ws.Range("E21").Value = "'0.66%"
$ws.Range("G21").Value = "'21"
$ws.Range("E22").Value = "'-1.82%"
$ws.Range("G22").Value = "'21"
$ws.Range("D23").Value = "'3.568"
$ws.Range("E23").Value = "'-4.94%"
$ws.Range("G23").Value = "'21"
$ws.Range("D24").Value = "'0.1374"
$ws.Range("E24").Value = "'1.76%"
$ws.Range("G24").Value = "'21"
$ws.Range("D25").Value = "'0.001215"
$ws.Range("E25").Value = "'-2.32%"
$ws.Range("G25").Value = "'21"
$ws.Range("D26").Value = "'0.004475"
$ws.Range("E26").Value = "'-2.36%"
$ws.Range("G26").Value = "'21"
$ws.Range("D27").Value = "'0.00009897"
$ws.Range("E27").Value = "'3.12%"
$ws.Range("G27").Value = "'21"
$ws.Range("E28").Value = "'-25.26%"
$ws.Range("G28").Value = "'21"
$ws.Range("G29").Value = "'21"
$ws.Range("G30").Value = "'21"
$ws.Range("G31").Value = "'21"
$ws.Range("G32").Value = "'21"
$ws.Range("G33").Value = "'21"
$ws.Range("G34").Value = "'21"
$ws.Range("G35").Value = "'21"
$ws.Range("G36").Value = "'21"
$ws.Range("G37").Value = "'21"
$ws.Range("G38").Value = "'21"
$ws.Range("G39").Value = "'21"
$ws.Range("D40").Value = "'0.03729"
$ws.Range("E40").Value = "'1.74%"
$ws.Range("G40").Value = "'21"
$ws.Range("D41").Value = "'0.005998"
$ws.Range("E41").Value = "'74.63%"
$ws.Range("G41").Value = "'21"
$ws.Range("E42").Value = "'-21.90%"
$ws.Range("G42").Value = "'21"
$ws.Range("D43").Value = "'0.002299"
$ws.Range("E43").Value = "'-13.54%"
$ws.Range("G43").Value = "'21"
$ws.Range("D44").Value = "'0.009712"
$ws.Range("E44").Value = "'17.71%"
$ws.Range("G44").Value = "'21"
$ws.Range("D45").Value = "'0.00005157"
$ws.Range("E45").Value = "'-4.31%"
$ws.Range("G45").Value = "'21"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("G46").Value = "'21"
$ws.Range("E47").Value = "'-8.21%"
$ws.Range("G47").Value = "'21"
$ws.Range("E48").Value = "'-3.27%"
$ws.Range("G48").Value = "'21"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("G49").Value = "'21"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("G50").Value = "'21"
$ws.Range("G51").Value = "'21"
